# Edit script: build up map_projected module changes reflected in KMeans.xlsx output
# - Add new column J ("clustering result") while column I becomes "TTP" (distance to nearest centroid)
# - Update a handful of SC(PPM) (column H) values that were recomputed upstream

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers: I1 becomes "TTP", J1 becomes "clustering result"
$ws.Cells.Item(1, 9).Value = "TTP"
$ws.Cells.Item(1, 10).Value = "clustering result"

# New "TTP" values for column I (rows 2-110)
$ttpValues = @(200.7860811985793,219.67328,146.8905782359963,194.3939531826629,139.5210013481891,141.269905274554,161.5074793128845,211.3318620951872,171.2410063241107,165.0336831099196,167.4266791044776,165.2836842105263,175.0338963446475,154.0116095238095,168.4999261682243,150.8921410404624,167.8427690763052,170.2389635627531,141.668238277512,166.823125490196,135.6634228855721,211.4708162758621,235.9211626329114,238.9125260491401,252.8184900585962,240.3466045629732,177.2628692307692,254.1618026273788,128.558256097561,176.5183923853211,178.3039076923077,201.5683957709251,189.2734051807229,177.2914351291513,192.5835129032258,127.9173560975609,209.8638640776699,175.2765445444422,185.7423717468094,263.0107908872902,146.1845155080214,136.2804266409266,231.2406951289398,160.7258966005666,126.7963636363637,222.8859755102041,133.6775954198473,190.6038227205312,218.349211721007,154.9662354727303,171.7672587813467,214.7173207171314,213.0718284518829,203.8279584659913,231.2100940170941,128.697895522388,186.3168966442953,209.0480375,192.7620289855073,163.9894539682539,163.628891902834,138.8930606060606,160.282187826087,201.7597771771772,197.6402116719243,180.1954190476191,192.1763793103448,192.8701104089219,181.8563975903615,152.7280308823529,164.3664111888112,170.1288148936171,134.0391735294118,132.2711643835617,128.5266777671411,138.5517161490683,122.0903581585345,157.1136561361632,119.9004795366795,178.7356301587301,188.2333185185186,147.2500602006689,166.6356631578947,117.8769655172414,260.0274157654389,118.1651999275037,122.8026928210982,113.134430395489,255.9015213004625,165.406010301293,168.693351334644,176.6202822889133,155.2334911659352,255.8861823529412,141.3051496240602,235.819891005291,255.2607464459498,135.3959147324642,178.9731962131889,168.0001698795181,123.86365,233.793043454039,136.2215,248.2423,230.2159946902655,157.2081818181818,266.0761746031746,211.7785130434783,134.6594236641221)

# "clustering result" values moved to column J (rows 2-110)
$clusterValues = @(6,5,0,5,0,0,2,5,2,2,2,2,2,2,2,2,2,2,0,2,0,5,1,1,1,1,2,1,0,2,2,5,4,2,5,0,4,4,4,1,0,0,1,2,0,3,0,5,5,2,2,5,5,5,1,0,5,4,4,2,2,0,2,4,4,4,4,4,2,2,2,2,0,0,7,0,0,2,0,2,5,0,2,0,1,0,0,7,1,2,2,2,2,1,0,1,1,0,2,2,7,1,0,1,1,2,1,5,0)

for ($i = 0; $i -lt $ttpValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $ttpValues[$i]
    $ws.Cells.Item($row, 10).Value = $clusterValues[$i]
}

# A few SC(PPM) (column H) values were recomputed upstream
$ws.Cells.Item(23, 8).Value = 70.80259014109356
$ws.Cells.Item(24, 8).Value = 70.80259014109356
$ws.Cells.Item(25, 8).Value = 70.80259014109356
$ws.Cells.Item(26, 8).Value = 70.80259014109356
$ws.Cells.Item(27, 8).Value = 70.80259014109356
$ws.Cells.Item(41, 8).Value = 70.80259014109356
